# Applies the "ERROR" (uncertainty) column addition + one extra (d,p) level row (L=3 at 2131 keV)
# split out of the old row 8, with every downstream row shifted / renumbered accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "ERROR" column header in E1; reuse the bordered/centered header style from D1 ---
$ws.Range("E1").Value = "ERROR"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 23 is brand new (table grew from 21 to 22 levels); give column A the table style ---
$ws.Range("A23").Value = 21
$ws.Range("A22").Copy()
$ws.Range("A23").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Final ENERGY/L/SPECTROSCOPIC_FACTOR/ERROR values for every data row, A2:E23 ---
$data = @(
    @(2, 0, 1862.87, 2, 0.05922545059314008, 0.00568035920186994),
    @(3, 1, 1875.2, 2, 0.0443288990743002, 0.003845078982516567),
    @(4, 2, 2054.045985291768, 2, 0.02810964105644485, 0.001491722245118135),
    @(5, 3, 2077.425494124893, 2, 0.0409853248420584, 0.001739056214628624),
    @(6, 4, 2095.238465959659, 2, 0.02100317822164607, 0.001324524752716419),
    @(7, 5, 2113.2, 2, 0.03969538251213013, 0.001665873366816012),
    @(8, 6, 2131.014162259718, 2, 0.006495147395202784, 0),
    @(9, 7, 2131.014162259718, 3, 0.01124859826003481, 0),
    @(10, 8, 2145.124556760028, 2, 0.008333543767399554, 0.0009152446912172919),
    @(11, 9, 2176.865287362391, 4, 0.07716494257580392, 0.008573882508422658),
    @(12, 10, 2196.568275413792, 2, 0.01350355374353616, 0.001129388131277569),
    @(13, 11, 2220.798153327344, 2, 0.01541937895772884, 0.001239499916216145),
    @(14, 12, 2246.353097309466, 1, 0.004273261164151795, 0.0002831192004141017),
    @(15, 13, 2306.889948659801, 2, 0.02531035941517577, 0.001752652179462273),
    @(16, 14, 2318.773858810782, 2, 0.04014875941325093, 0.001971127723131183),
    @(17, 15, 2332.404601796536, 3, 0.009930304388552249, 0.001429963831951524),
    @(18, 16, 2342.559621553097, 0, 0.008687218081936166, 0.001116928039106079),
    @(19, 17, 2352.409431043697, 2, 0.01189349736647337, 0.001427219683976804),
    @(20, 18, 2364.788262037358, 2, 0.01937268493144125, 0.001596649856986916),
    @(21, 19, 2374.896726123765, 2, 0.04401072620265029, 0.002088100148300926),
    @(22, 20, 2384.426696065908, 1, 0.002501150571464353, 0.0002881559463663034),
    @(23, 21, 2397.177155842639, 3, 0.05384121018720114, 0.002575736907103759)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
}

$ws.Range("A1").Select() | Out-Null
